# Remove the "blog" paragraph from the Teaching Statement's Evaluations
# section (the paragraph beginning "With regard to assisting (primarily
# graduate) students ..." and ending "...or a similar question."),
# including its paragraph mark, leaving "Evaluations." directly followed
# by "Teaching Awards".

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute(
    "With regard to assisting*or a similar question.",
    $false, $false, $true, $false, $false,
    $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the target paragraph to delete."
}

# $rng now spans exactly the found text (real Word Find.Execute behavior).
# Grab the whole paragraph (text + trailing paragraph mark) via the
# paragraph's own Range so the paragraph itself is removed, not just its
# text, and delete it.
$para = $rng.Paragraphs(1)
$paraRange = $para.Range
$paraRange.Delete()

Write-Output "Deleted paragraph. Remaining paragraph count: $($d.Paragraphs.Count)"
